# v2.0: Ranking dinamico por grupos SISBEN y nuevo shapefile
# Correccion de duplicados en datos (LA FLORA, PARQUE ENTRENUBES)
#
# The source table had two duplicate UPZ rows that need to be removed:
#   - row 99  -> UPZ "LA FLORA" / Localidad "Usme" (duplicate entry)
#   - row 115 -> UPZ "PARQUE ENTRENUBES" / Localidad "Usme" (duplicate, all-zero row)
#
# Removing these two rows shifts every following row up, which also keeps the
# RANKING column (column A, already stored as row-number-1) consistent without
# any further edits, and shrinks the used range from A1:R115 down to A1:R113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "LA FLORA" duplicate row first.
$ws.Rows(99).Delete()

# After the row above is removed, everything shifted up by one, so the former
# row 115 ("PARQUE ENTRENUBES" duplicate) is now row 114. Delete it too.
$ws.Rows(114).Delete()
